# Updated symbol list on Fri Dec 16 14:26:58 UTC 2022 with GitHub Actions
# Refreshes the "Price" column with newer quotes and swaps the
# GateToken / KuCoinToken rows (6 and 7), including their updated prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 / Row 7: GateToken and KuCoinToken swap places ---
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("E6").Value = "5GateTokenGT"

$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("E7").Value = "6KuCoinTokenKCS"

# --- Price ("D") column refresh ---
# The Price column stores numeric-looking values as plain text (inline
# strings) in the workbook, so each value is written with a leading
# apostrophe (forces Excel to keep it as text instead of silently
# converting it to a number) and the cell style is then reset back to
# "Normal" so no stray number-format / quote-prefix styling lingers on
# the cell.
$prices = @{
    "D2"  = "249.39"
    "D3"  = "24.32"
    "D4"  = "5.926"
    "D5"  = "0.05891"
    "D6"  = "3.428"
    "D7"  = "6.576"
    "D8"  = "1.335"
    "D9"  = "0.7971"
    "D10" = "0.1480"
    "D11" = "0.07733"
    "D12" = "0.03323"
    "D13" = "0.03027"
    "D14" = "0.09254"
    "D15" = "3.556"
    "D16" = "0.001668"
    "D17" = "0.04754"
    "D18" = "0.0006041"
    "D19" = "0.006273"
    "D20" = "0.005552"
    "D21" = "0.001065"
    "D22" = "0.0001497"
    "D23" = "3.726"
    "D25" = "0.3355"
    "D26" = "0.1256"
    "D27" = "0.0006488"
    "D40" = "0.04383"
    "D41" = "0.007038"
    "D42" = "0.1065"
    "D43" = "0.003366"
    "D44" = "0.01008"
    "D45" = "0.002465"
    "D46" = "0.00005902"
    "D47" = "0.00000000752"
    "D48" = "0.9921"
    "D49" = "0.1106"
    "D50" = "0.00002104"
    "D51" = "0.01012"
}

foreach ($ref in $prices.Keys) {
    $ws.Range($ref).Value = "'" + $prices[$ref]
    $ws.Range($ref).Style = "Normal"
}
